# Auto-generated from diff: updates market-price/profit columns (H-N)
# across all 8 job sheets to reflect refreshed Universalis pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 526.96155
$ws.Range("I33").Value = 168.8
$ws.Range("J33").Value = 1015.36365
$ws.Range("K33").Value = 168.8
$ws.Range("L33").Value = 1015.36365
$ws.Range("M33").Value = 60.19999999999999
$ws.Range("N33").Value = -1473.36365

$ws.Range("H94").Value = 2360.4
$ws.Range("I94").Value = 2360.4
$ws.Range("K94").Value = 2360.4
$ws.Range("M94").Value = -1909.4

$ws.Range("H96").Value = 244.33333
$ws.Range("I96").Value = 252.41667
$ws.Range("J96").Value = 228.16667
$ws.Range("K96").Value = 757.25001
$ws.Range("L96").Value = 684.50001
$ws.Range("M96").Value = 615.74999
$ws.Range("N96").Value = -3430.50001

$ws.Range("H97").Value = 1637.5
$ws.Range("J97").Value = 1936.6666
$ws.Range("L97").Value = 5809.9998
$ws.Range("N97").Value = -6801.9998

$ws.Range("H99").Value = 859.05554
$ws.Range("I99").Value = 357.5
$ws.Range("J99").Value = 1260.3
$ws.Range("K99").Value = 1072.5
$ws.Range("L99").Value = 3780.9
$ws.Range("M99").Value = 425.5
$ws.Range("N99").Value = -6776.9

$ws.Range("H104").Value = 941.5
$ws.Range("I104").Value = 941.5
$ws.Range("K104").Value = 2824.5
$ws.Range("M104").Value = -1077.5

$ws.Range("H112").Value = 43479840
$ws.Range("J112").Value = 1817
$ws.Range("L112").Value = 5451
$ws.Range("N112").Value = -7667

$ws.Range("H132").Value = 711688.75
$ws.Range("I132").Value = 1325.6774
$ws.Range("J132").Value = 7003475.5
$ws.Range("K132").Value = 3977.0322
$ws.Range("L132").Value = 21010426.5
$ws.Range("M132").Value = -1447.0322
$ws.Range("N132").Value = -21015486.5

$ws.Range("H138").Value = 3474238.8
$ws.Range("I138").Value = 1612.862
$ws.Range("J138").Value = 8774563
$ws.Range("K138").Value = 4838.586
$ws.Range("L138").Value = 26323689
$ws.Range("M138").Value = 301.4139999999998
$ws.Range("N138").Value = -26333969

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2902.1
$ws.Range("I2").Value = 3017.4285
$ws.Range("J2").Value = 2633
$ws.Range("K2").Value = 3017.4285
$ws.Range("L2").Value = 2633
$ws.Range("M2").Value = -2904.4285
$ws.Range("N2").Value = -2859

$ws.Range("H32").Value = 1025.64
$ws.Range("I32").Value = 871.4699000000001
$ws.Range("J32").Value = 1778.3529
$ws.Range("K32").Value = 871.4699000000001
$ws.Range("L32").Value = 1778.3529
$ws.Range("M32").Value = -584.4699000000001
$ws.Range("N32").Value = -2352.3529

$ws.Range("H45").Value = 1433.2778
$ws.Range("I45").Value = 945
$ws.Range("K45").Value = 945
$ws.Range("M45").Value = -568

$ws.Range("H61").Value = 22773566
$ws.Range("I61").Value = 28600646
$ws.Range("J61").Value = 112703.11
$ws.Range("K61").Value = 28600646
$ws.Range("L61").Value = 112703.11
$ws.Range("M61").Value = -28600434
$ws.Range("N61").Value = -113127.11

$ws.Range("H74").Value = 7413464.5
$ws.Range("I74").Value = 10041431
$ws.Range("J74").Value = 113557.11
$ws.Range("K74").Value = 10041431
$ws.Range("L74").Value = 113557.11
$ws.Range("M74").Value = -10040557
$ws.Range("N74").Value = -115305.11

$ws.Range("H77").Value = 7413464.5
$ws.Range("I77").Value = 10041431
$ws.Range("J77").Value = 113557.11
$ws.Range("K77").Value = 50207155
$ws.Range("L77").Value = 567785.55
$ws.Range("M77").Value = -50202787
$ws.Range("N77").Value = -576521.55

$ws.Range("H97").Value = 2718341
$ws.Range("I97").Value = 4167629.5
$ws.Range("K97").Value = 4167629.5
$ws.Range("M97").Value = -4167133.5

$ws.Range("H116").Value = 2902.1
$ws.Range("I116").Value = 3017.4285
$ws.Range("J116").Value = 2633
$ws.Range("K116").Value = 3017.4285
$ws.Range("L116").Value = 2633
$ws.Range("M116").Value = -723.4285
$ws.Range("N116").Value = -7221

$ws.Range("H136").Value = 22773566
$ws.Range("I136").Value = 28600646
$ws.Range("J136").Value = 112703.11
$ws.Range("K136").Value = 85801938
$ws.Range("L136").Value = 338109.33
$ws.Range("M136").Value = -85799388
$ws.Range("N136").Value = -343209.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2902.1
$ws.Range("I3").Value = 3017.4285
$ws.Range("J3").Value = 2633
$ws.Range("K3").Value = 3017.4285
$ws.Range("L3").Value = 2633
$ws.Range("M3").Value = -2903.4285
$ws.Range("N3").Value = -2861

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H80").Value = 638.7857
$ws.Range("I80").Value = 309.42856
$ws.Range("K80").Value = 309.42856
$ws.Range("M80").Value = 688.5714399999999

$ws.Range("H83").Value = 638.7857
$ws.Range("I83").Value = 309.42856
$ws.Range("K83").Value = 1547.1428
$ws.Range("M83").Value = 3444.8572

$ws.Range("H94").Value = 887.55554
$ws.Range("I94").Value = 839.6
$ws.Range("J94").Value = 947.5
$ws.Range("K94").Value = 839.6
$ws.Range("L94").Value = 947.5
$ws.Range("M94").Value = -388.6
$ws.Range("N94").Value = -1849.5

$ws.Range("H99").Value = 1446.6666
$ws.Range("I99").Value = 1407.1428
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1407.1428
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 90.85719999999992
$ws.Range("N99").Value = -4996

$ws.Range("H134").Value = 2322.7
$ws.Range("I134").Value = 1604.6923
$ws.Range("J134").Value = 2871.7646
$ws.Range("K134").Value = 4814.0769
$ws.Range("L134").Value = 8615.293799999999
$ws.Range("M134").Value = -2279.0769
$ws.Range("N134").Value = -13685.2938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1637
$ws.Range("I16").Value = 1455.5
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1455.5
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1168.5
$ws.Range("N16").Value = -2574

$ws.Range("H22").Value = 250000660
$ws.Range("I22").Value = 1000000000
$ws.Range("J22").Value = 865.3333
$ws.Range("K22").Value = 1000000000
$ws.Range("L22").Value = 865.3333
$ws.Range("M22").Value = -999999650
$ws.Range("N22").Value = -1565.3333

$ws.Range("H31").Value = 1780.9207
$ws.Range("I31").Value = 1056.9778
$ws.Range("J31").Value = 3590.7778
$ws.Range("K31").Value = 1056.9778
$ws.Range("L31").Value = 3590.7778
$ws.Range("M31").Value = -761.9777999999999
$ws.Range("N31").Value = -4180.7778

$ws.Range("H34").Value = 1780.9207
$ws.Range("I34").Value = 1056.9778
$ws.Range("J34").Value = 3590.7778
$ws.Range("K34").Value = 1056.9778
$ws.Range("L34").Value = 3590.7778
$ws.Range("M34").Value = -854.9777999999999
$ws.Range("N34").Value = -3994.7778

$ws.Range("H113").Value = 1637
$ws.Range("I113").Value = 1455.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1455.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 714.5
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1032.8983
$ws.Range("I131").Value = 421.2857
$ws.Range("J131").Value = 1223.1777
$ws.Range("K131").Value = 1263.8571
$ws.Range("L131").Value = 3669.5331
$ws.Range("M131").Value = 3776.1429
$ws.Range("N131").Value = -13749.5331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2311.111
$ws.Range("I97").Value = 2362.5
$ws.Range("K97").Value = 2362.5
$ws.Range("M97").Value = -1866.5

$ws.Range("H132").Value = 78785.766
$ws.Range("I132").Value = 51551.55
$ws.Range("J132").Value = 169566.5
$ws.Range("K132").Value = 154654.65
$ws.Range("L132").Value = 508699.5
$ws.Range("M132").Value = -152124.65
$ws.Range("N132").Value = -513759.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20788.588
$ws.Range("I132").Value = 10134.147
$ws.Range("J132").Value = 54994.95
$ws.Range("K132").Value = 30402.441
$ws.Range("L132").Value = 164984.85
$ws.Range("M132").Value = -27872.441
$ws.Range("N132").Value = -170044.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2356.7144
$ws.Range("I96").Value = 2166.3333
$ws.Range("J96").Value = 2499.5
$ws.Range("K96").Value = 2166.3333
$ws.Range("L96").Value = 2499.5
$ws.Range("M96").Value = -793.3332999999998
$ws.Range("N96").Value = -5245.5

$ws.Range("H100").Value = 53590.42
$ws.Range("I100").Value = 56003.332
$ws.Range("K100").Value = 112006.664
$ws.Range("M100").Value = -111465.664
